# Zeitplan.xlsx update: mark several tasks complete, replace the "Abgabe in"
# countdown formula with a typed date, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")
$ws.Activate()

# "Abgabe in" cell: previously =TODAY()+27, now a manually typed date text.
$ws.Range("E6").Value = "13.12."

# Mark the following tasks' "% erledigt" as 100% done (the dependent
# "Fortschritt" column recalculates automatically from the table formula).
$ws.Range("E11").Value = 1   # Storyboard-Zeichungen
$ws.Range("E12").Value = 1   # Tonspur für Storyboard-Video
$ws.Range("E14").Value = 1   # Storyboard-Video erstellen
$ws.Range("E17").Value = 1   # Auto ins Git committen
$ws.Range("E22").Value = 1   # Holztexturen -  (Bett, Schrank, Nachttisch)

# Move the visible selection/scroll position like the author left it.
[void]$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
